$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.29%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.40%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.597"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.15%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.02%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.43%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.370"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.95%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.896"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.73%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.817"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.76%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9434"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.77%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1192"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.51%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1917"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.55%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09890"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.34%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04341"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "11.77%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.85%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001273"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.96%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005946"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.72%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.531"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.77%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3536"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.07%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.736"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "8.98%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1369"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.16%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2496"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.44%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04386"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.62%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001244"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.99%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004343"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.69%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001235"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.75%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004005"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "31.53%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02784"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.38%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05728"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.95%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007945"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.35%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.47%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1422"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.96%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002104"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.64%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01006"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.39%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007316"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.34%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.25%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003449"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.89%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002279"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002109"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.25%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002008"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.25%"
